$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.4214906692504883
$ws.Range("E2").Value = 391.5080634653295
$ws.Range("F2").Value = 0.01487401281013962
$ws.Range("G2").Value = 0.01173913137544259
$ws.Range("H2").Value = 0.01068783509471432
$ws.Range("I2").Value = 0.009330681130351386
$ws.Range("J2").Value = 0.009141809519973794
$ws.Range("K2").Value = 0.008686094586635155
$ws.Range("L2").Value = 0.008543903460618079
$ws.Range("M2").Value = 0.008439079480885601
$ws.Range("N2").Value = 0.008439079480885601
$ws.Range("O2").Value = 0.008219063424769423
$ws.Range("P2").Value = 0.008143757360676334
$ws.Range("Q2").Value = 0.008055718503233117
$ws.Range("R2").Value = 0.007790026963101831
$ws.Range("S2").Value = 0.007790026963101831
$ws.Range("T2").Value = 0.007770669148991911
$ws.Range("U2").Value = 0.007762232450338706
$ws.Range("V2").Value = 0.007694904104483942
$ws.Range("W2").Value = 0.007690447973004655
$ws.Range("X2").Value = 0.007673125165817578
$ws.Range("Y2").Value = 0.007631736129928449

$ws.Range("C3").Value = 0.4155347347259521
$ws.Range("E3").Value = 394.4358961441503
$ws.Range("F3").Value = 0.01486086559180789
$ws.Range("G3").Value = 0.01239571431091643
$ws.Range("H3").Value = 0.01089659202358738
$ws.Range("I3").Value = 0.01041545646954961
$ws.Range("J3").Value = 0.009876703316521029
$ws.Range("K3").Value = 0.009383835173560257
$ws.Range("L3").Value = 0.009085341256392855
$ws.Range("M3").Value = 0.008949734960971968
$ws.Range("N3").Value = 0.008749440036196517
$ws.Range("O3").Value = 0.008590248980310088
$ws.Range("P3").Value = 0.00809356801349421
$ws.Range("Q3").Value = 0.008081664771651198
$ws.Range("R3").Value = 0.008037154035295736
$ws.Range("S3").Value = 0.007900937955367617
$ws.Range("T3").Value = 0.00783369743418815
$ws.Range("U3").Value = 0.007824327946359878
$ws.Range("V3").Value = 0.007700937777668978
$ws.Range("W3").Value = 0.007700937777668978
$ws.Range("X3").Value = 0.007700937777668978
$ws.Range("Y3").Value = 0.007688808891698836

$ws.Range("C4").Value = 0.3906524181365967
$ws.Range("E4").Value = 393.3327093337593
$ws.Range("F4").Value = 0.01501075095763927
$ws.Range("G4").Value = 0.01225556827196887
$ws.Range("H4").Value = 0.01066016398700701
$ws.Range("I4").Value = 0.009719307002354551
$ws.Range("J4").Value = 0.009436522884056806
$ws.Range("K4").Value = 0.009110746838417569
$ws.Range("L4").Value = 0.008776658150022553
$ws.Range("M4").Value = 0.008637853744198523
$ws.Range("N4").Value = 0.008637853744198523
$ws.Range("O4").Value = 0.008509734735444653
$ws.Range("P4").Value = 0.008454456300108796
$ws.Range("Q4").Value = 0.008366361152612301
$ws.Range("R4").Value = 0.008117919846180027
$ws.Range("S4").Value = 0.008117919846180027
$ws.Range("T4").Value = 0.00807236740263502
$ws.Range("U4").Value = 0.007946264902342053
$ws.Range("V4").Value = 0.007851891865966409
$ws.Range("W4").Value = 0.007766691473560014
$ws.Range("X4").Value = 0.007713018742021399
$ws.Range("Y4").Value = 0.007667304275511877

$ws.Range("C5").Value = 0.3749721050262451
$ws.Range("E5").Value = 374.0917530445331
$ws.Range("F5").Value = 0.0156223292812564
$ws.Range("G5").Value = 0.01164470385329851
$ws.Range("H5").Value = 0.01046936646909238
$ws.Range("I5").Value = 0.009604811391567897
$ws.Range("J5").Value = 0.009221860737736639
$ws.Range("K5").Value = 0.008746951909401578
$ws.Range("L5").Value = 0.008438361041362884
$ws.Range("M5").Value = 0.008408543188531382
$ws.Range("N5").Value = 0.008242790314507854
$ws.Range("O5").Value = 0.007916860302503994
$ws.Range("P5").Value = 0.007805632465968913
$ws.Range("Q5").Value = 0.007647903673662484
$ws.Range("R5").Value = 0.007613875656365545
$ws.Range("S5").Value = 0.007613875656365545
$ws.Range("T5").Value = 0.007398628357417084
$ws.Range("U5").Value = 0.007368133604985428
$ws.Range("V5").Value = 0.007347951223551188
$ws.Range("W5").Value = 0.007292236901452885
$ws.Range("X5").Value = 0.007292236901452885
$ws.Range("Y5").Value = 0.007292236901452885

$ws.Range("C6").Value = 0.4062495231628418
$ws.Range("E6").Value = 395.2947960853926
$ws.Range("F6").Value = 0.01535112010910288
$ws.Range("G6").Value = 0.01217937147587525
$ws.Range("H6").Value = 0.01052218124013061
$ws.Range("I6").Value = 0.01009822869905259
$ws.Range("J6").Value = 0.009590056364027651
$ws.Range("K6").Value = 0.009158424328888229
$ws.Range("L6").Value = 0.009010402458605358
$ws.Range("M6").Value = 0.008964887933980772
$ws.Range("N6").Value = 0.008285141132711114
$ws.Range("O6").Value = 0.008285141132711114
$ws.Range("P6").Value = 0.008285141132711114
$ws.Range("Q6").Value = 0.008285141132711114
$ws.Range("R6").Value = 0.008079737399571863
$ws.Range("S6").Value = 0.007925033604594865
$ws.Range("T6").Value = 0.007803733034794261
$ws.Range("U6").Value = 0.007795111721787047
$ws.Range("V6").Value = 0.007795111721787047
$ws.Range("W6").Value = 0.007795111721787047
$ws.Range("X6").Value = 0.007738850780379154
$ws.Range("Y6").Value = 0.007705551580611939

$ws.Range("C7").Value = 0.3750381469726562
$ws.Range("E7").Value = 378.7706983889238
$ws.Range("F7").Value = 0.01484058267703947
$ws.Range("G7").Value = 0.01219741434529057
$ws.Range("H7").Value = 0.01081476701275156
$ws.Range("I7").Value = 0.009959888420104198
$ws.Range("J7").Value = 0.008926599690856421
$ws.Range("K7").Value = 0.008721063625342397
$ws.Range("L7").Value = 0.00858125004455544
$ws.Range("M7").Value = 0.00823856093919377
$ws.Range("N7").Value = 0.008060018359920347
$ws.Range("O7").Value = 0.008000562139575592
$ws.Range("P7").Value = 0.00793307105900572
$ws.Range("Q7").Value = 0.007811222961919983
$ws.Range("R7").Value = 0.007639663146828378
$ws.Range("S7").Value = 0.007639663146828378
$ws.Range("T7").Value = 0.007552981895750009
$ws.Range("U7").Value = 0.007537117415090069
$ws.Range("V7").Value = 0.007447566121157855
$ws.Range("W7").Value = 0.007447472835343896
$ws.Range("X7").Value = 0.007383444413039449
$ws.Range("Y7").Value = 0.007383444413039449

$ws.Range("C8").Value = 0.3749895095825195
$ws.Range("E8").Value = 404.5853246511815
$ws.Range("F8").Value = 0.01501684414691057
$ws.Range("G8").Value = 0.01294771677502004
$ws.Range("H8").Value = 0.01124731311590158
$ws.Range("I8").Value = 0.01017757595711667
$ws.Range("J8").Value = 0.009376488273568819
$ws.Range("K8").Value = 0.009266777950411515
$ws.Range("L8").Value = 0.009266777950411515
$ws.Range("M8").Value = 0.008922273908277109
$ws.Range("N8").Value = 0.00871258672989635
$ws.Range("O8").Value = 0.008655507376042343
$ws.Range("P8").Value = 0.008384059003399679
$ws.Range("Q8").Value = 0.008296645302643756
$ws.Range("R8").Value = 0.008228209099681447
$ws.Range("S8").Value = 0.00815822404091217
$ws.Range("T8").Value = 0.008056274562204209
$ws.Range("U8").Value = 0.008029837251642444
$ws.Range("V8").Value = 0.007963885250049246
$ws.Range("W8").Value = 0.00791049510503755
$ws.Range("X8").Value = 0.00791049510503755
$ws.Range("Y8").Value = 0.007886653501972346

$ws.Range("C9").Value = 0.46872878074646
$ws.Range("E9").Value = 384.9595354773101
$ws.Range("F9").Value = 0.01508986923990492
$ws.Range("G9").Value = 0.01281312083344518
$ws.Range("H9").Value = 0.01071001888151002
$ws.Range("I9").Value = 0.01010054734706933
$ws.Range("J9").Value = 0.01000299496701331
$ws.Range("K9").Value = 0.009414665097331441
$ws.Range("L9").Value = 0.00868647390723868
$ws.Range("M9").Value = 0.008625935912295686
$ws.Range("N9").Value = 0.008413005382741036
$ws.Range("O9").Value = 0.008111643211449256
$ws.Range("P9").Value = 0.008018997626378003
$ws.Range("Q9").Value = 0.00779837097444869
$ws.Range("R9").Value = 0.00779837097444869
$ws.Range("S9").Value = 0.007725747153539908
$ws.Range("T9").Value = 0.007555142356729959
$ws.Range("U9").Value = 0.007553997729293456
$ws.Range("V9").Value = 0.007509587242226616
$ws.Range("W9").Value = 0.007509587242226616
$ws.Range("X9").Value = 0.007509587242226616
$ws.Range("Y9").Value = 0.007504084512228264

$ws.Range("C10").Value = 0.4487597942352295
$ws.Range("E10").Value = 374.6258951282925
$ws.Range("F10").Value = 0.0154082046939301
$ws.Range("G10").Value = 0.01226765608829645
$ws.Range("H10").Value = 0.01059133391809279
$ws.Range("I10").Value = 0.009852641974054159
$ws.Range("J10").Value = 0.009209363183251579
$ws.Range("K10").Value = 0.008787487375102215
$ws.Range("L10").Value = 0.008428239056580425
$ws.Range("M10").Value = 0.00837327334979724
$ws.Range("N10").Value = 0.008169419065630987
$ws.Range("O10").Value = 0.007993284622946083
$ws.Range("P10").Value = 0.007874369314754376
$ws.Range("Q10").Value = 0.007726887298274774
$ws.Range("R10").Value = 0.007634650047909592
$ws.Range("S10").Value = 0.007572541165067709
$ws.Range("T10").Value = 0.007566519234626902
$ws.Range("U10").Value = 0.007425684563618809
$ws.Range("V10").Value = 0.007425684563618809
$ws.Range("W10").Value = 0.007337526725766675
$ws.Range("X10").Value = 0.007337526725766675
$ws.Range("Y10").Value = 0.007302649027841958

$ws.Range("C11").Value = 0.4062039852142334
$ws.Range("E11").Value = 402.1185997150969
$ws.Range("F11").Value = 0.01391244208316284
$ws.Range("G11").Value = 0.01172603935901344
$ws.Range("H11").Value = 0.01100454952501474
$ws.Range("I11").Value = 0.01093676855389992
$ws.Range("J11").Value = 0.01003773686803902
$ws.Range("K11").Value = 0.009717608735873866
$ws.Range("L11").Value = 0.009282803829237537
$ws.Range("M11").Value = 0.00907272228495651
$ws.Range("N11").Value = 0.008784075319131144
$ws.Range("O11").Value = 0.008611235810706833
$ws.Range("P11").Value = 0.008459596386065563
$ws.Range("Q11").Value = 0.008309884376095885
$ws.Range("R11").Value = 0.008256626452859126
$ws.Range("S11").Value = 0.008138812950416185
$ws.Range("T11").Value = 0.008028112115027847
$ws.Range("U11").Value = 0.007986934398802589
$ws.Range("V11").Value = 0.007967765787042662
$ws.Range("W11").Value = 0.007934488099060409
$ws.Range("X11").Value = 0.007934488099060409
$ws.Range("Y11").Value = 0.007838569195226059
